$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7423.385
$ws.Range("I74").Value = 14000
$ws.Range("J74").Value = 5450.4
$ws.Range("K74").Value = 14000
$ws.Range("L74").Value = 5450.4
$ws.Range("M74").Value = -13064
$ws.Range("N74").Value = -7322.4

$ws.Range("H77").Value = 7423.385
$ws.Range("I77").Value = 14000
$ws.Range("J77").Value = 5450.4
$ws.Range("K77").Value = 70000
$ws.Range("L77").Value = 27252
$ws.Range("M77").Value = -65320
$ws.Range("N77").Value = -36612

$ws.Range("H113").Value = 3197.1
$ws.Range("J113").Value = 3096
$ws.Range("L113").Value = 3096
$ws.Range("N113").Value = -9604

$ws.Range("H133").Value = 52614.285
$ws.Range("J133").Value = 52614.285
$ws.Range("L133").Value = 52614.285
$ws.Range("N133").Value = -62734.285

$ws.Range("H135").Value = 6336.75
$ws.Range("I135").Value = 6270.8
$ws.Range("J135").Value = 6366.727
$ws.Range("K135").Value = 56437.2
$ws.Range("L135").Value = 57300.543
$ws.Range("M135").Value = -53902.2
$ws.Range("N135").Value = -62370.543

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1326
$ws.Range("I74").Value = 1171.091
$ws.Range("J74").Value = 1468
$ws.Range("K74").Value = 1171.091
$ws.Range("L74").Value = 1468
$ws.Range("M74").Value = -297.0909999999999
$ws.Range("N74").Value = -3216

$ws.Range("H77").Value = 1326
$ws.Range("I77").Value = 1171.091
$ws.Range("J77").Value = 1468
$ws.Range("K77").Value = 5855.455
$ws.Range("L77").Value = 7340
$ws.Range("M77").Value = -1487.455
$ws.Range("N77").Value = -16076

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 11136.363
$ws.Range("J35").Value = 11136.363
$ws.Range("L35").Value = 11136.363
$ws.Range("N35").Value = -11756.363

$ws.Range("H86").Value = 55613560
$ws.Range("I86").Value = 90911640
$ws.Range("J86").Value = 145145.86
$ws.Range("K86").Value = 90911640
$ws.Range("L86").Value = 145145.86
$ws.Range("M86").Value = -90910517
$ws.Range("N86").Value = -147391.86

$ws.Range("H89").Value = 55613560
$ws.Range("I89").Value = 90911640
$ws.Range("J89").Value = 145145.86
$ws.Range("K89").Value = 454558200
$ws.Range("L89").Value = 725729.2999999999
$ws.Range("M89").Value = -454552584
$ws.Range("N89").Value = -736961.2999999999

$ws.Range("H105").Value = 9182.200000000001
$ws.Range("I105").Value = 8800
$ws.Range("K105").Value = 8800
$ws.Range("M105").Value = -7053

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1239.8462
$ws.Range("I16").Value = 1217.2858
$ws.Range("J16").Value = 1266.1666
$ws.Range("K16").Value = 1217.2858
$ws.Range("L16").Value = 1266.1666
$ws.Range("M16").Value = -930.2858000000001
$ws.Range("N16").Value = -1840.1666

$ws.Range("H50").Value = 6092
$ws.Range("J50").Value = 6092
$ws.Range("L50").Value = 6092
$ws.Range("N50").Value = -7342

$ws.Range("H51").Value = 20000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 20000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 20000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -21472

$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 20000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 20000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -20696

$ws.Range("H86").Value = 71846.62
$ws.Range("I86").Value = 18399.666
$ws.Range("J86").Value = 117658.29
$ws.Range("K86").Value = 18399.666
$ws.Range("L86").Value = 117658.29
$ws.Range("M86").Value = -17276.666
$ws.Range("N86").Value = -119904.29

$ws.Range("H89").Value = 71846.62
$ws.Range("I89").Value = 18399.666
$ws.Range("J89").Value = 117658.29
$ws.Range("K89").Value = 91998.33
$ws.Range("L89").Value = 588291.45
$ws.Range("M89").Value = -86382.33
$ws.Range("N89").Value = -599523.45

$ws.Range("H105").Value = 796.8
$ws.Range("I105").Value = 605.2174
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 605.2174
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = 1141.7826
$ws.Range("N105").Value = -6494

$ws.Range("H113").Value = 1239.8462
$ws.Range("I113").Value = 1217.2858
$ws.Range("J113").Value = 1266.1666
$ws.Range("K113").Value = 1217.2858
$ws.Range("L113").Value = 1266.1666
$ws.Range("M113").Value = 952.7141999999999
$ws.Range("N113").Value = -5606.1666

$ws.Range("H132").Value = 46056.824
$ws.Range("I132").Value = 2238.125
$ws.Range("J132").Value = 146213.86
$ws.Range("K132").Value = 6714.375
$ws.Range("L132").Value = 438641.58
$ws.Range("M132").Value = -4184.375
$ws.Range("N132").Value = -443701.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 4105.4546
$ws.Range("J43").Value = 5993.3335
$ws.Range("L43").Value = 17980.0005
$ws.Range("N43").Value = -18208.0005

$ws.Range("H98").Value = 4500.5264
$ws.Range("J98").Value = 5214.2666
$ws.Range("L98").Value = 15642.7998
$ws.Range("N98").Value = -18638.7998

$ws.Range("H140").Value = 4784.326
$ws.Range("I140").Value = 3148.5
$ws.Range("J140").Value = 5500
$ws.Range("K140").Value = 9445.5
$ws.Range("L140").Value = 16500
$ws.Range("M140").Value = -4265.5
$ws.Range("N140").Value = -26860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26792614
$ws.Range("I70").Value = 56254960
$ws.Range("J70").Value = 8663.637000000001
$ws.Range("K70").Value = 56254960
$ws.Range("L70").Value = 8663.637000000001
$ws.Range("M70").Value = -56254690
$ws.Range("N70").Value = -9203.637000000001

$ws.Range("H73").Value = 26792614
$ws.Range("I73").Value = 56254960
$ws.Range("J73").Value = 8663.637000000001
$ws.Range("K73").Value = 56254960
$ws.Range("L73").Value = 8663.637000000001
$ws.Range("M73").Value = -56254024
$ws.Range("N73").Value = -10535.637

$ws.Range("H126").Value = 16677664
$ws.Range("I126").Value = 15210
$ws.Range("K126").Value = 45630
$ws.Range("M126").Value = -43160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 457340.8
$ws.Range("I132").Value = 716335.4
$ws.Range("J132").Value = 4100.375
$ws.Range("K132").Value = 2149006.2
$ws.Range("L132").Value = 12301.125
$ws.Range("M132").Value = -2146476.2
$ws.Range("N132").Value = -17361.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 37427.375
$ws.Range("I46").Value = 22000
$ws.Range("J46").Value = 39631.285
$ws.Range("K46").Value = 22000
$ws.Range("L46").Value = 39631.285
$ws.Range("M46").Value = -21769
$ws.Range("N46").Value = -40093.285

$ws.Range("H132").Value = 60763790
$ws.Range("I132").Value = 187501650
$ws.Range("J132").Value = 2269391.5
$ws.Range("K132").Value = 562504950
$ws.Range("L132").Value = 6808174.5
$ws.Range("M132").Value = -562502420
$ws.Range("N132").Value = -6813234.5

$ws.Range("H133").Value = 44095
$ws.Range("I133").Value = 39000
$ws.Range("K133").Value = 39000
$ws.Range("M133").Value = -33940

$ws.Range("H134").Value = 37427.375
$ws.Range("I134").Value = 22000
$ws.Range("J134").Value = 39631.285
$ws.Range("K134").Value = 66000
$ws.Range("L134").Value = 118893.855
$ws.Range("M134").Value = -63465
$ws.Range("N134").Value = -123963.855
